# ES2N-Requisitos Funcionais: add RF15 / RF16 rows, clean up stray
# "Pode executar RF13 TIRAR" note and an accidental underline on B12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the leftover underline formatting on B12 ("Gerenciar
#    fornecedores") so it matches the rest of the Nome column.
$ws.Range("A12").Copy()
$ws.Range("B12").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# 2) Clear the stray "Pode executar RF13 TIRAR" note in E17 and restore
#    the plain (non-highlighted) look used by the rest of column E.
$ws.Range("E10").Copy()
$ws.Range("E17").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("E17").ClearContents()

# 3) Turn the two placeholder rows (24/25) into real RF15 / RF16 entries.
#    First drop the yellow highlight from B24/B25 so they read like normal
#    "Nome" cells.
$ws.Range("A18").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("B25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A18").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C18").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("C25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$excel.CutCopyMode = $false

# Fill in the actual requirement data (codes first, then descriptions,
# then priorities - matches the order new shared strings were appended).
$ws.Range("A24").Value = "RF15"
$ws.Range("A25").Value = "RF16"

$ws.Range("C24").Value = "Por meio de filtros o gerente e o laboratório podem consultar laudos de matéria-prima."
$ws.Range("C25").Value = "Por meio de filtros o gerente e o laboratório podem consultar laudos de produto."

$ws.Range("D24").Value = "Média"
$ws.Range("D25").Value = "Média"
